# new analyse result ,bad!!
$wb = $excel.ActiveWorkbook

# --- Sheet "5m6s" (sheet1): update row 5, add row 6 ---
$ws = $wb.Worksheets.Item("5m6s")
$ws.Range("C5").Value = 1.216
$ws.Range("D5").Value = 21.98
$ws.Range("E5").Value = "0/12.5/87.5"
$ws.Range("F5").Value = 0.2143
$ws.Range("G5").Value = 0
$ws.Range("G5").NumberFormat = "0%"
$ws.Range("H5").Value = -5.848

$ws.Range("A6").Value = "geo_mini"
$ws.Range("B6").Value = 0.001
$ws.Range("C6").Value = 0.237
$ws.Range("D6").Value = 4.4
$ws.Range("E6").Value = "4.17/2.08/93.75"
$ws.Range("F6").Value = 0
$ws.Range("F6").NumberFormat = "0.00%"
$ws.Range("G6").Value = 0
$ws.Range("G6").NumberFormat = "0%"
$ws.Range("H6").Value = -2.081
$ws.Range("H6").Select()

# --- Sheet "5nbz" (sheet2): update row 5, add row 6 ---
$ws = $wb.Worksheets.Item("5nbz")
$ws.Range("B5").Value = 0.016
$ws.Range("C5").Value = 2.796
$ws.Range("D5").Value = 65.98
$ws.Range("E5").Value = "2.98/25.11/71.91"
$ws.Range("F5").Value = 0.1058
$ws.Range("G5").Value = 0.0131
$ws.Range("H5").Value = -5.608

$ws.Range("A6").Value = "geo_mini"
$ws.Range("B6").Value = 0.001
$ws.Range("C6").Value = 0.304
$ws.Range("D6").Value = 0.93
$ws.Range("E6").Value = "0.78/5.99/93.23"
$ws.Range("F6").Value = 0.52
$ws.Range("F6").NumberFormat = "0.00%"
$ws.Range("G6").Value = 0
$ws.Range("G6").NumberFormat = "0.00%"
$ws.Range("H6").Value = -2.754
$ws.Range("L4").Select()

# --- Sheet "6h3l" (sheet3): add row 6 ---
$ws = $wb.Worksheets.Item("6h3l")
$ws.Range("A6").Value = "geo_mini"
$ws.Range("B6").Value = 0.001
$ws.Range("C6").Value = 0.326
$ws.Range("D6").Value = 3.04
$ws.Range("E6").Value = "0.43/11.34/88.23"
$ws.Range("F6").Value = 0.0013
$ws.Range("F6").NumberFormat = "0.00%"
$ws.Range("G6").Value = 0
$ws.Range("G6").NumberFormat = "0.00%"
$ws.Range("H6").Value = -3.585
$ws.Range("H6").Select()

# --- Sheet "6h3n" (sheet4): add row 6 ---
$ws = $wb.Worksheets.Item("6h3n")
$ws.Range("A6").Value = "geo_mini"
$ws.Range("B6").Value = 0.001
$ws.Range("C6").Value = 0.33
$ws.Range("D6").Value = 2.29
$ws.Range("E6").Value = "0.27/9.90/89.93"
$ws.Range("F6").Value = 0
$ws.Range("F6").NumberFormat = "0.00%"
$ws.Range("G6").Value = 0
$ws.Range("G6").NumberFormat = "0%"
$ws.Range("H6").Value = -3.347
$ws.Range("E11").Select()

# --- Sheet "6n2p" (sheet5): add row 6 ---
$ws = $wb.Worksheets.Item("6n2p")
$ws.Range("A6").Value = "geo_mini"
$ws.Range("B6").Value = 0.001
$ws.Range("C6").Value = 0.356
$ws.Range("D6").Value = 3.57
$ws.Range("E6").Value = "1.16/13.26/85.58"
$ws.Range("F6").Value = 0.012
$ws.Range("F6").NumberFormat = "0.00%"
$ws.Range("G6").Value = 0
$ws.Range("G6").NumberFormat = "0%"
$ws.Range("H6").Value = -4.11
$ws.Range("H6").Select()

# --- Sheet "6o1k" (sheet6): add row 6 ---
$ws = $wb.Worksheets.Item("6o1k")
$ws.Range("A6").Value = "geo_mini"
$ws.Range("B6").Value = 0.001
$ws.Range("C6").Value = 0.333
$ws.Range("D6").Value = 3.88
$ws.Range("E6").Value = "1.24/13.11/85.65"
$ws.Range("F6").Value = 0.0042
$ws.Range("F6").NumberFormat = "0.00%"
$ws.Range("G6").Value = 0
$ws.Range("G6").NumberFormat = "0%"
$ws.Range("H6").Value = -3.618
$ws.Range("H6").Select()

# --- Sheet "6o1l" (sheet7): add row 6 ---
$ws = $wb.Worksheets.Item("6o1l")
$ws.Range("A6").Value = "geo_mini"
$ws.Range("B6").Value = 0.001
$ws.Range("C6").Value = 0.331
$ws.Range("D6").Value = 3.69
$ws.Range("E6").Value = "2.07/13.36/84.57"
$ws.Range("F6").Value = 0.0049
$ws.Range("F6").NumberFormat = "0.00%"
$ws.Range("G6").Value = 0
$ws.Range("G6").NumberFormat = "0%"
$ws.Range("H6").Value = -3.891
$ws.Range("H13").Select()

# --- Sheet "6o1m" (sheet8): add row 6 (G6 left default, no % style) ---
$ws = $wb.Worksheets.Item("6o1m")
$ws.Range("A6").Value = "geo_mini"
$ws.Range("B6").Value = 0.001
$ws.Range("C6").Value = 0.332
$ws.Range("D6").Value = 4.35
$ws.Range("E6").Value = "1.82/12.66/85.52"
$ws.Range("F6").Value = 0.0061
$ws.Range("F6").NumberFormat = "0.00%"
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = -3.793
$ws.Range("H15").Select()
